$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.415.70"
$ws.Range("E2").Value = "  +6.31%  "

$ws.Range("D3").Value = "'1.814.79"
$ws.Range("E3").Value = "  +6.29%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'344.92"
$ws.Range("E5").Value = "  +4.10%  "

$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.3847"
$ws.Range("E7").Value = "  +4.54%  "

$ws.Range("D8").Value = "'50.35"
$ws.Range("E8").Value = "  +4.07%  "

$ws.Range("D9").Value = "'0.3527"
$ws.Range("E9").Value = "  +6.75%  "

$ws.Range("D10").Value = "'1.240"
$ws.Range("E10").Value = "  +5.96%  "

$ws.Range("E11").Value = "  +5.81%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").Value = "'22.63"
$ws.Range("E13").Value = "  +13.44%  "

$ws.Range("D14").Value = "'6.646"
$ws.Range("E14").Value = "  +7.17%  "

$ws.Range("D15").Value = "'7.246"
$ws.Range("E15").Value = "  +5.59%  "

$ws.Range("D16").Value = "'1.813.85"
$ws.Range("E16").Value = "  +6.78%  "

$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "  +4.84%  "

$ws.Range("D18").Value = "'0.06799"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("D19").Value = "'87.16"
$ws.Range("E19").Value = "  +7.17%  "

$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "'17.87"
$ws.Range("E21").Value = "  +10.37%  "

$ws.Range("D22").Value = "'6.554"
$ws.Range("E22").Value = "  +8.05%  "

$ws.Range("E23").Value = "  +1.55%  "

$ws.Range("D24").Value = "'27.413.54"
$ws.Range("E24").Value = "  +6.44%  "

$ws.Range("D25").Value = "'2.472"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").Value = "'2.728"
$ws.Range("E26").Value = "  +9.75%  "

$ws.Range("D27").Value = "'22.42"
$ws.Range("E27").Value = "  +17.02%  "

$ws.Range("D28").Value = "'1.500"
$ws.Range("E28").Value = "  +15.71%  "

$ws.Range("D29").Value = "'154.47"
$ws.Range("E29").Value = "  +3.22%  "

$ws.Range("D30").Value = "'2.020.11"
$ws.Range("E30").Value = "  +6.98%  "

$ws.Range("D31").Value = "'137.00"
$ws.Range("E31").Value = "  +6.97%  "

$ws.Range("D32").Value = "'6.406"
$ws.Range("E32").Value = "  +7.56%  "

$ws.Range("D33").Value = "'4.103"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").Value = "'13.88"
$ws.Range("E34").Value = "  +7.74%  "

$ws.Range("D35").Value = "'0.08831"
$ws.Range("E35").Value = "  +3.91%  "

$ws.Range("D36").Value = "'1.722"
$ws.Range("E36").Value = "  +2.81%  "

$ws.Range("D37").Value = "'5.658"
$ws.Range("E37").Value = "  +6.34%  "

$ws.Range("D38").Value = "'0.7099"
$ws.Range("E38").Value = "  +16.02%  "

$ws.Range("D39").Value = "'0.06569"
$ws.Range("E39").Value = "  +5.59%  "

$ws.Range("D40").Value = "'0.2274"
$ws.Range("E40").Value = "  +7.02%  "

$ws.Range("D41").Value = "'0.02418"
$ws.Range("E41").Value = "  +7.04%  "

$ws.Range("D42").Value = "'9.031"
$ws.Range("E42").Value = "  +5.67%  "

$ws.Range("D43").Value = "'1.260"
$ws.Range("E43").Value = "  -1.05%  "

$ws.Range("D44").Value = "'15.00"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("D45").Value = "'0.6637"
$ws.Range("E45").Value = "  +13.69%  "

$ws.Range("D46").Value = "'0.9993"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "'3.980"

$ws.Range("D48").Value = "'2.196"
$ws.Range("E48").Value = "  +9.55%  "

$ws.Range("D49").Value = "'133.01"
$ws.Range("E49").Value = "  +5.26%  "

$ws.Range("D50").Value = "'0.07366"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("D51").Value = "'80.89"
$ws.Range("E51").Value = "  +5.68%  "
